# Adding HUB test case for Language check
#
# - Removes the stale "testing123" placeholder cells (B5/B6) from the "hub" sheet
# - Moves the hub sheet's active selection to F11
# - Adds a new "system" worksheet (after "testdata") that reuses the hub
#   sheet's header-row formatting, with a new "Title"/"systeem" data row

$wb  = $excel.ActiveWorkbook
$hub = $wb.Worksheets.Item("hub")

# Clear the old test placeholder values (previously "testing123") from B5:B6
$null = $hub.Range("B5").ClearContents()
$null = $hub.Range("B6").ClearContents()

# Add the new "system" sheet right after the last existing sheet ("testdata")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sys = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$sys.Name = "system"

# Reuse hub's styled header row (A1:M1) for the new sheet, formatting included
$hub.Range("A1:M1").Copy($sys.Range("A1"))

# New data row
$sys.Range("A2").Value = "Title"
$sys.Range("B2").Value = "systeem"

# Selection on the new sheet
$null = $sys.Range("A3").Select()

# Restore "hub" as the active sheet and set its new selection
$null = $hub.Activate()
$null = $hub.Range("F11").Select()
